$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The simulation was rerun: two new entries ("Holden", "Rizzie Spiral") were
# inserted into the material list right after "Spiral5", "Thomas Hex" was
# renamed to "Matthies Hex", and the per-material results table (rows 4-31)
# was recomputed/shifted down accordingly, adding two new trailing rows.

# Give the two newly appended index rows (30 and 31) the same formatting
# (bold, centered, bordered) as column A of the existing rows.
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 4
$ws.Cells.Item(4, 2).Value2 = 'Holden'
$ws.Cells.Item(4, 3).Value2 = 1.001871610715191
$ws.Cells.Item(4, 4).Value2 = 1.000701857040393
$ws.Cells.Item(4, 5).Value2 = 0.997192584478829
$ws.Cells.Item(4, 6).Value2 = 0.997192584478829
$ws.Cells.Item(4, 7).Value2 = 0.9994385162656434
$ws.Cells.Item(4, 8).Value2 = 0.9981746025113138
$ws.Cells.Item(4, 9).Value2 = 1.000701857040393
$ws.Cells.Item(4, 10).Value2 = 0.9993967543059022
$ws.Cells.Item(4, 11).Value2 = 1.001042090195771
$ws.Cells.Item(4, 12).Value2 = 1.001871610715191
$ws.Cells.Item(4, 13).Value2 = 1.001871610715191
$ws.Cells.Item(4, 14).Value2 = 1.001871610715191
$ws.Cells.Item(4, 15).Value2 = 1.000701857040393
$ws.Cells.Item(4, 16).Value2 = 0.9989472207596108
$ws.Cells.Item(4, 17).Value2 = 1.000049305673147
$ws.Cells.Item(4, 18).Value2 = 0.9999220174114711
$ws.Cells.Item(4, 19).Value2 = 0.9990970652750413
$ws.Cells.Item(4, 20).Value2 = 0.9999220174114711
$ws.Cells.Item(4, 21).Value2 = 0.9997907016350789
$ws.Cells.Item(4, 22).Value2 = 1.000206883451102
$ws.Cells.Item(4, 23).Value2 = 0.9998149840691797
# Row 5
$ws.Cells.Item(5, 2).Value2 = 'Rizzie Spiral'
$ws.Cells.Item(5, 3).Value2 = 1.003657505162014
$ws.Cells.Item(5, 4).Value2 = 1.001371574925625
$ws.Cells.Item(5, 5).Value2 = 0.994513752543125
$ws.Cells.Item(5, 6).Value2 = 0.994513752543125
$ws.Cells.Item(5, 7).Value2 = 0.9989027489004727
$ws.Cells.Item(5, 8).Value2 = 0.9964328069176354
$ws.Cells.Item(5, 9).Value2 = 1.001371574925625
$ws.Cells.Item(5, 10).Value2 = 0.9988211298422112
$ws.Cells.Item(5, 11).Value2 = 1.002036455582161
$ws.Cells.Item(5, 12).Value2 = 1.003657505162014
$ws.Cells.Item(5, 13).Value2 = 1.003657505162014
$ws.Cells.Item(5, 14).Value2 = 1.003657505162014
$ws.Cells.Item(5, 15).Value2 = 1.001371574925625
$ws.Cells.Item(5, 16).Value2 = 0.9979426637343749
$ws.Cells.Item(5, 17).Value2 = 1.000096352383918
$ws.Cells.Item(5, 18).Value2 = 0.9998476108769211
$ws.Cells.Item(5, 19).Value2 = 0.9982354857703203
$ws.Cells.Item(5, 20).Value2 = 0.9998476108769211
$ws.Cells.Item(5, 21).Value2 = 0.9995909906182436
$ws.Cells.Item(5, 22).Value2 = 1.000404293526998
$ws.Cells.Item(5, 23).Value2 = 0.9996384435998585
# Row 6
$ws.Cells.Item(6, 2).Value2 = 'RotRing OmegaMax-90'
$ws.Cells.Item(6, 3).Value2 = 1.001057162099064
$ws.Cells.Item(6, 4).Value2 = 1.000396441248276
$ws.Cells.Item(6, 5).Value2 = 0.9984142526515971
$ws.Cells.Item(6, 6).Value2 = 0.9984142526515971
$ws.Cells.Item(6, 7).Value2 = 0.9996828502631481
$ws.Cells.Item(6, 8).Value2 = 0.9989689382146675
$ws.Cells.Item(6, 9).Value2 = 1.000396441248276
$ws.Cells.Item(6, 10).Value2 = 0.9996592623850588
$ws.Cells.Item(6, 11).Value2 = 1.000588620668503
$ws.Cells.Item(6, 12).Value2 = 1.001057162099064
$ws.Cells.Item(6, 13).Value2 = 1.001057162099064
$ws.Cells.Item(6, 14).Value2 = 1.001057162099064
$ws.Cells.Item(6, 15).Value2 = 1.000396441248276
$ws.Cells.Item(6, 16).Value2 = 0.9994053469499368
$ws.Cells.Item(6, 17).Value2 = 1.000027851816668
$ws.Cells.Item(6, 18).Value2 = 0.9999559519996458
$ws.Cells.Item(6, 19).Value2 = 0.9994899854283107
$ws.Cells.Item(6, 20).Value2 = 0.9999559519996458
$ws.Cells.Item(6, 21).Value2 = 0.9998817795959991
$ws.Cells.Item(6, 22).Value2 = 1.000116856096612
$ws.Cells.Item(6, 23).Value2 = 0.999895496097324
# Row 7
$ws.Cells.Item(7, 2).Value2 = 'Equal Angle'
$ws.Cells.Item(7, 3).Value2 = 1.001965226174351
$ws.Cells.Item(7, 4).Value2 = 1.000736962716139
$ws.Cells.Item(7, 5).Value2 = 0.9970521631988478
$ws.Cells.Item(7, 6).Value2 = 0.9970521631988478
$ws.Cells.Item(7, 7).Value2 = 0.9994104319092214
$ws.Cells.Item(7, 8).Value2 = 0.9980832996253617
$ws.Cells.Item(7, 9).Value2 = 1.000736962716139
$ws.Cells.Item(7, 10).Value2 = 0.9993665807492768
$ws.Cells.Item(7, 11).Value2 = 1.001094214423631
$ws.Cells.Item(7, 12).Value2 = 1.001965226174351
$ws.Cells.Item(7, 13).Value2 = 1.001965226174351
$ws.Cells.Item(7, 14).Value2 = 1.001965226174351
$ws.Cells.Item(7, 15).Value2 = 1.000736962716139
$ws.Cells.Item(7, 16).Value2 = 0.9988945629574935
$ws.Cells.Item(7, 17).Value2 = 1.000051771732708
$ws.Cells.Item(7, 18).Value2 = 0.9999181173631126
$ws.Cells.Item(7, 19).Value2 = 0.9990519022214213
$ws.Cells.Item(7, 20).Value2 = 0.9999181173631126
$ws.Cells.Item(7, 21).Value2 = 0.9997802332096537
$ws.Cells.Item(7, 22).Value2 = 1.000217231802593
$ws.Cells.Item(7, 23).Value2 = 0.999805730189121
# Row 8
$ws.Cells.Item(8, 2).Value2 = 'Tilt Rotate'
$ws.Cells.Item(8, 3).Value2 = 1.006536151903368
$ws.Cells.Item(8, 4).Value2 = 1.002451052123913
$ws.Cells.Item(8, 5).Value2 = 0.9901958024679519
$ws.Cells.Item(8, 6).Value2 = 0.9901958024679519
$ws.Cells.Item(8, 7).Value2 = 0.9980391592460238
$ws.Cells.Item(8, 8).Value2 = 0.9936252546539249
$ws.Cells.Item(8, 9).Value2 = 1.002451052123913
$ws.Cells.Item(8, 10).Value2 = 0.9978933136236334
$ws.Cells.Item(8, 11).Value2 = 1.003639249598744
$ws.Cells.Item(8, 12).Value2 = 1.006536151903368
$ws.Cells.Item(8, 13).Value2 = 1.006536151903368
$ws.Cells.Item(8, 14).Value2 = 1.006536151903368
$ws.Cells.Item(8, 15).Value2 = 1.002451052123913
$ws.Cells.Item(8, 16).Value2 = 0.9963234272959325
$ws.Cells.Item(8, 17).Value2 = 1.000172182873773
$ws.Cells.Item(8, 18).Value2 = 0.9997276688317442
$ws.Cells.Item(8, 19).Value2 = 0.9968467227384995
$ws.Cells.Item(8, 20).Value2 = 0.9997276688317444
$ws.Cells.Item(8, 21).Value2 = 0.9992690800297166
$ws.Cells.Item(8, 22).Value2 = 1.000722494404447
$ws.Cells.Item(8, 23).Value2 = 0.9993538794676839
# Row 9
$ws.Cells.Item(9, 2).Value2 = 'CLR'
$ws.Cells.Item(9, 3).Value2 = 1.000216580719323
$ws.Cells.Item(9, 4).Value2 = 1.000081222221127
$ws.Cells.Item(9, 5).Value2 = 0.9996751277528185
$ws.Cells.Item(9, 6).Value2 = 0.9996751277528185
$ws.Cells.Item(9, 7).Value2 = 0.9999350254851215
$ws.Cells.Item(9, 8).Value2 = 0.9997887664999148
$ws.Cells.Item(9, 9).Value2 = 1.000081222221127
$ws.Cells.Item(9, 10).Value2 = 0.9999301932228801
$ws.Cells.Item(9, 11).Value2 = 1.000120591847271
$ws.Cells.Item(9, 12).Value2 = 1.000216580719323
$ws.Cells.Item(9, 13).Value2 = 1.000216580719323
$ws.Cells.Item(9, 14).Value2 = 1.000216580719323
$ws.Cells.Item(9, 15).Value2 = 1.000081222221127
$ws.Cells.Item(9, 16).Value2 = 0.9998781749869725
$ws.Cells.Item(9, 17).Value2 = 1.000005707722003
$ws.Cells.Item(9, 18).Value2 = 0.999990976897756
$ws.Cells.Item(9, 19).Value2 = 0.9998955143989416
$ws.Cells.Item(9, 20).Value2 = 0.9999909768977561
$ws.Cells.Item(9, 21).Value2 = 0.9999757809790371
$ws.Cells.Item(9, 22).Value2 = 1.000023940927094
$ws.Cells.Item(9, 23).Value2 = 0.9999785912461978
# Row 10
$ws.Cells.Item(10, 2).Value2 = 'Rizzie Hex'
$ws.Cells.Item(10, 3).Value2 = 1.000016284515201
$ws.Cells.Item(10, 4).Value2 = 1.000006110863086
$ws.Cells.Item(10, 5).Value2 = 0.9999755728388378
$ws.Cells.Item(10, 6).Value2 = 0.9999755728388378
$ws.Cells.Item(10, 7).Value2 = 0.9999951147354407
$ws.Cells.Item(10, 8).Value2 = 0.9999841174356137
$ws.Cells.Item(10, 9).Value2 = 1.000006110863086
$ws.Cells.Item(10, 10).Value2 = 0.9999947517045974
$ws.Cells.Item(10, 11).Value2 = 1.000009069157193
$ws.Cells.Item(10, 12).Value2 = 1.000016284515201
$ws.Cells.Item(10, 13).Value2 = 1.000016284515201
$ws.Cells.Item(10, 14).Value2 = 1.000016284515201
$ws.Cells.Item(10, 15).Value2 = 1.000006110863086
$ws.Cells.Item(10, 16).Value2 = 0.9999908418509618
$ws.Cells.Item(10, 17).Value2 = 1.000000431283842
$ws.Cells.Item(10, 18).Value2 = 0.9999993227390416
$ws.Cells.Item(10, 19).Value2 = 0.9999921451355069
$ws.Cells.Item(10, 20).Value2 = 0.9999993227390416
$ws.Cells.Item(10, 21).Value2 = 0.9999981799804305
$ws.Cells.Item(10, 22).Value2 = 1.000001800887385
$ws.Cells.Item(10, 23).Value2 = 0.9999983915141319
# Row 11
$ws.Cells.Item(11, 2).Value2 = 'Matthies Hex'
$ws.Cells.Item(11, 3).Value2 = 1.000359495831011
$ws.Cells.Item(11, 4).Value2 = 1.000134817007111
$ws.Cells.Item(11, 5).Value2 = 0.9994607561731895
$ws.Cells.Item(11, 6).Value2 = 0.9994607561731895
$ws.Cells.Item(11, 7).Value2 = 0.9998921514297636
$ws.Cells.Item(11, 8).Value2 = 0.9996493815667122
$ws.Cells.Item(11, 9).Value2 = 1.000134817007111
$ws.Cells.Item(11, 10).Value2 = 0.9998841291328039
$ws.Cells.Item(11, 11).Value2 = 1.000200166106374
$ws.Cells.Item(11, 12).Value2 = 1.000359495831011
$ws.Cells.Item(11, 13).Value2 = 1.000359495831011
$ws.Cells.Item(11, 14).Value2 = 1.000359495831011
$ws.Cells.Item(11, 15).Value2 = 1.000134817007111
$ws.Cells.Item(11, 16).Value2 = 0.9997977865901504
$ws.Cells.Item(11, 17).Value2 = 1.000009473069958
$ws.Cells.Item(11, 18).Value2 = 0.9999850230037707
$ws.Cells.Item(11, 19).Value2 = 0.9998265674377015
$ws.Cells.Item(11, 20).Value2 = 0.9999850230037707
$ws.Cells.Item(11, 21).Value2 = 0.999959799536029
$ws.Cells.Item(11, 22).Value2 = 1.000039738795026
$ws.Cells.Item(11, 23).Value2 = 0.9999644642817596
# Row 12
$ws.Cells.Item(12, 2).Value2 = 'Tilt Rotate_Partial'
$ws.Cells.Item(12, 3).Value2 = 1.006633772398039
$ws.Cells.Item(12, 4).Value2 = 1.00248766002167
$ws.Cells.Item(12, 5).Value2 = 0.9900493723235453
$ws.Cells.Item(12, 6).Value2 = 0.9900493723235453
$ws.Cells.Item(12, 7).Value2 = 0.9980098726036252
$ws.Cells.Item(12, 8).Value2 = 0.993530045268928
$ws.Cells.Item(12, 9).Value2 = 1.00248766002167
$ws.Cells.Item(12, 10).Value2 = 0.9978618484041302
$ws.Cells.Item(12, 11).Value2 = 1.003693602984971
$ws.Cells.Item(12, 12).Value2 = 1.006633772398039
$ws.Cells.Item(12, 13).Value2 = 1.006633772398039
$ws.Cells.Item(12, 14).Value2 = 1.006633772398039
$ws.Cells.Item(12, 15).Value2 = 1.00248766002167
$ws.Cells.Item(12, 16).Value2 = 0.9962685161726076
$ws.Cells.Item(12, 17).Value2 = 1.0001747542129
$ws.Cells.Item(12, 18).Value2 = 0.9997236015810845
$ws.Cells.Item(12, 19).Value2 = 0.9967996269164484
$ws.Cells.Item(12, 20).Value2 = 0.9997236015810845
$ws.Cells.Item(12, 21).Value2 = 0.999258163286846
$ws.Cells.Item(12, 22).Value2 = 1.000733285109085
$ws.Cells.Item(12, 23).Value2 = 0.9993442292533222
# Row 13
$ws.Cells.Item(13, 2).Value2 = 'RotRing OmegaMax-60'
$ws.Cells.Item(13, 3).Value2 = 1.007127898833531
$ws.Cells.Item(13, 4).Value2 = 1.002672968615168
$ws.Cells.Item(13, 5).Value2 = 0.9893081457438621
$ws.Cells.Item(13, 6).Value2 = 0.9893081457438621
$ws.Cells.Item(13, 7).Value2 = 0.9978616300395569
$ws.Cells.Item(13, 8).Value2 = 0.9930480958903197
$ws.Cells.Item(13, 9).Value2 = 1.002672968615168
$ws.Cells.Item(13, 10).Value2 = 0.9977025781810431
$ws.Cells.Item(13, 11).Value2 = 1.00396872645385
$ws.Cells.Item(13, 12).Value2 = 1.007127898833531
$ws.Cells.Item(13, 13).Value2 = 1.007127898833531
$ws.Cells.Item(13, 14).Value2 = 1.007127898833531
$ws.Cells.Item(13, 15).Value2 = 1.002672968615168
$ws.Cells.Item(13, 16).Value2 = 0.9959905571795149
$ws.Cells.Item(13, 17).Value2 = 1.000187773398105
$ws.Cells.Item(13, 18).Value2 = 0.9997030043975205
$ws.Cells.Item(13, 19).Value2 = 0.9965612308466909
$ws.Cells.Item(13, 20).Value2 = 0.9997030043975205
$ws.Cells.Item(13, 21).Value2 = 0.9992028978434011
$ws.Cells.Item(13, 22).Value2 = 1.000787898041427
$ws.Cells.Item(13, 23).Value2 = 0.9992953765465624
# Row 14
$ws.Cells.Item(14, 2).Value2 = 'Equal Angle_Partial'
$ws.Cells.Item(14, 3).Value2 = 1.002229103852631
$ws.Cells.Item(14, 4).Value2 = 1.000835916273683
$ws.Cells.Item(14, 5).Value2 = 0.9966563488210515
$ws.Cells.Item(14, 6).Value2 = 0.9966563488210515
$ws.Cells.Item(14, 7).Value2 = 0.9993312680526325
$ws.Cells.Item(14, 8).Value2 = 0.997825938821054
$ws.Cells.Item(14, 9).Value2 = 1.000835916273683
$ws.Cells.Item(14, 10).Value2 = 0.99928152896842
$ws.Cells.Item(14, 11).Value2 = 1.001241137431578
$ws.Cells.Item(14, 12).Value2 = 1.002229103852631
$ws.Cells.Item(14, 13).Value2 = 1.002229103852631
$ws.Cells.Item(14, 14).Value2 = 1.002229103852631
$ws.Cells.Item(14, 15).Value2 = 1.000835916273683
$ws.Cells.Item(14, 16).Value2 = 0.9987461325473672
$ws.Cells.Item(14, 17).Value2 = 1.000058722621052
$ws.Cells.Item(14, 18).Value2 = 0.9999071229824552
$ws.Cells.Item(14, 19).Value2 = 0.9989245980210515
$ws.Cells.Item(14, 20).Value2 = 0.9999071229824552
$ws.Cells.Item(14, 21).Value2 = 0.9997507244789464
$ws.Cells.Item(14, 22).Value2 = 1.000246400353683
$ws.Cells.Item(14, 23).Value2 = 0.9997796448118417
# Row 15
$ws.Cells.Item(15, 2).Value2 = 'Rizzie Hex_Partial'
$ws.Cells.Item(15, 3).Value2 = 0.9986953291366857
$ws.Cells.Item(15, 4).Value2 = 0.9995107532491248
$ws.Cells.Item(15, 5).Value2 = 1.001957007022278
$ws.Cells.Item(15, 6).Value2 = 1.001957007022278
$ws.Cells.Item(15, 7).Value2 = 1.000391400449006
$ws.Cells.Item(15, 8).Value2 = 1.001272457802474
$ws.Cells.Item(15, 9).Value2 = 0.9995107532491248
$ws.Cells.Item(15, 10).Value2 = 1.000420511660059
$ws.Cells.Item(15, 11).Value2 = 0.9992735782526495
$ws.Cells.Item(15, 12).Value2 = 0.9986953291366857
$ws.Cells.Item(15, 13).Value2 = 0.9986953291366857
$ws.Cells.Item(15, 14).Value2 = 0.9986953291366857
$ws.Cells.Item(15, 15).Value2 = 0.9995107532491248
$ws.Cells.Item(15, 16).Value2 = 1.000733880135701
$ws.Cells.Item(15, 17).Value2 = 0.9999656324545917
$ws.Cells.Item(15, 18).Value2 = 1.000054363136029
$ws.Cells.Item(15, 19).Value2 = 1.000629423977154
$ws.Cells.Item(15, 20).Value2 = 1.000054363136029
$ws.Cells.Item(15, 21).Value2 = 1.000145900267037
$ws.Cells.Item(15, 22).Value2 = 0.9998557860409665
$ws.Cells.Item(15, 23).Value2 = 1.000128973852675
# Row 16
$ws.Cells.Item(16, 2).Value2 = 'ND Single'
$ws.Cells.Item(16, 3).Value2 = 1.011461099999999
$ws.Cells.Item(16, 4).Value2 = 1.004297899999999
$ws.Cells.Item(16, 9).Value2 = 1.004297899999999
$ws.Cells.Item(16, 12).Value2 = 1.011461099999999
$ws.Cells.Item(16, 13).Value2 = 1.011461099999999
$ws.Cells.Item(16, 14).Value2 = 1.011461099999999
$ws.Cells.Item(16, 15).Value2 = 1.004297899999999
$ws.Cells.Item(16, 17).Value2 = 1.000301919999999
# Row 17
$ws.Cells.Item(17, 2).Value2 = 'RD Single'
$ws.Cells.Item(17, 3).Value2 = 0.95415576
$ws.Cells.Item(17, 4).Value2 = 0.98280841
$ws.Cells.Item(17, 5).Value2 = 1.0687664
$ws.Cells.Item(17, 6).Value2 = 1.0687664
$ws.Cells.Item(17, 7).Value2 = 1.0137533
$ws.Cells.Item(17, 8).Value2 = 1.0447123
$ws.Cells.Item(17, 9).Value2 = 0.98280841
$ws.Cells.Item(17, 10).Value2 = 1.0147762
$ws.Cells.Item(17, 11).Value2 = 0.9744745399999999
$ws.Cells.Item(17, 12).Value2 = 0.95415576
$ws.Cells.Item(17, 13).Value2 = 0.95415576
$ws.Cells.Item(17, 14).Value2 = 0.95415576
$ws.Cells.Item(17, 15).Value2 = 0.98280841
$ws.Cells.Item(17, 16).Value2 = 1.025787405
$ws.Cells.Item(17, 17).Value2 = 0.998792305
$ws.Cells.Item(17, 18).Value2 = 1.00191019
$ws.Cells.Item(17, 19).Value2 = 1.022117003333333
$ws.Cells.Item(17, 20).Value2 = 1.00191019
$ws.Cells.Item(17, 21).Value2 = 1.0051266925
$ws.Cells.Item(17, 22).Value2 = 0.9949325059999999
$ws.Cells.Item(17, 23).Value2 = 1.004531915
# Row 18
$ws.Cells.Item(18, 2).Value2 = 'TD Single'
$ws.Cells.Item(18, 3).Value2 = 1.0114611
$ws.Cells.Item(18, 4).Value2 = 1.0042979
$ws.Cells.Item(18, 5).Value2 = 0.98280841
$ws.Cells.Item(18, 6).Value2 = 0.98280841
$ws.Cells.Item(18, 7).Value2 = 0.9965616799999999
$ws.Cells.Item(18, 8).Value2 = 0.98882193
$ws.Cells.Item(18, 9).Value2 = 1.0042979
$ws.Cells.Item(18, 10).Value2 = 0.9963059399999998
$ws.Cells.Item(18, 11).Value2 = 1.0063814
$ws.Cells.Item(18, 12).Value2 = 1.0114611
$ws.Cells.Item(18, 13).Value2 = 1.0114611
$ws.Cells.Item(18, 14).Value2 = 1.0114611
$ws.Cells.Item(18, 15).Value2 = 1.0042979
$ws.Cells.Item(18, 16).Value2 = 0.9935531550000001
$ws.Cells.Item(18, 17).Value2 = 1.00030192
$ws.Cells.Item(18, 18).Value2 = 0.9995224700000001
$ws.Cells.Item(18, 19).Value2 = 0.99447075
$ws.Cells.Item(18, 20).Value2 = 0.9995224700000001
$ws.Cells.Item(18, 21).Value2 = 0.9987183374999999
$ws.Cells.Item(18, 22).Value2 = 1.00126689
$ws.Cells.Item(18, 23).Value2 = 0.9988670325
# Row 19
$ws.Cells.Item(19, 2).Value2 = 'Morris Single'
$ws.Cells.Item(19, 3).Value2 = 0.96848208
$ws.Cells.Item(19, 4).Value2 = 0.98818078
$ws.Cells.Item(19, 5).Value2 = 1.0472769
$ws.Cells.Item(19, 6).Value2 = 1.0472769
$ws.Cells.Item(19, 7).Value2 = 1.0094554
$ws.Cells.Item(19, 8).Value2 = 1.0307397
$ws.Cells.Item(19, 9).Value2 = 0.98818078
$ws.Cells.Item(19, 10).Value2 = 1.0101587
$ws.Cells.Item(19, 11).Value2 = 0.9824512399999999
$ws.Cells.Item(19, 12).Value2 = 0.96848208
$ws.Cells.Item(19, 13).Value2 = 0.96848208
$ws.Cells.Item(19, 14).Value2 = 0.96848208
$ws.Cells.Item(19, 15).Value2 = 0.98818078
$ws.Cells.Item(19, 16).Value2 = 1.01772884
$ws.Cells.Item(19, 17).Value2 = 0.99916974
$ws.Cells.Item(19, 18).Value2 = 1.001313253333333
$ws.Cells.Item(19, 19).Value2 = 1.01520546
$ws.Cells.Item(19, 20).Value2 = 1.001313253333333
$ws.Cells.Item(19, 21).Value2 = 1.003524615
$ws.Cells.Item(19, 22).Value2 = 0.9965161080000001
$ws.Cells.Item(19, 23).Value2 = 1.0031156975
# Row 20
$ws.Cells.Item(20, 2).Value2 = 'Ring Perpendicular to ND'
$ws.Cells.Item(20, 3).Value2 = 1.003611022465753
$ws.Cells.Item(20, 4).Value2 = 1.001354144520548
$ws.Cells.Item(20, 5).Value2 = 0.9945834764383562
$ws.Cells.Item(20, 6).Value2 = 0.9945834764383562
$ws.Cells.Item(20, 7).Value2 = 0.9989166939726029
$ws.Cells.Item(20, 8).Value2 = 0.9964781419178085
$ws.Cells.Item(20, 9).Value2 = 1.001354144520548
$ws.Cells.Item(20, 10).Value2 = 0.9988361117808215
$ws.Cells.Item(20, 11).Value2 = 1.002010575616438
$ws.Cells.Item(20, 12).Value2 = 1.003611022465753
$ws.Cells.Item(20, 13).Value2 = 1.003611022465753
$ws.Cells.Item(20, 14).Value2 = 1.003611022465753
$ws.Cells.Item(20, 15).Value2 = 1.001354144520548
$ws.Cells.Item(20, 16).Value2 = 0.997968810479452
$ws.Cells.Item(20, 17).Value2 = 1.000095128150685
$ws.Cells.Item(20, 18).Value2 = 0.9998495478082191
$ws.Cells.Item(20, 19).Value2 = 0.9982579109132418
$ws.Cells.Item(20, 20).Value2 = 0.9998495478082191
$ws.Cells.Item(20, 21).Value2 = 0.9995961888013698
$ws.Cells.Item(20, 22).Value2 = 1.000399155534246
$ws.Cells.Item(20, 23).Value2 = 0.9996430389041095
# Row 21
$ws.Cells.Item(21, 2).Value2 = 'Ring Perpendicular to RD'
$ws.Cells.Item(21, 3).Value2 = 0.984316445263158
$ws.Cells.Item(21, 4).Value2 = 0.9941186657894736
$ws.Cells.Item(21, 5).Value2 = 1.023525342105263
$ws.Cells.Item(21, 6).Value2 = 1.023525342105263
$ws.Cells.Item(21, 7).Value2 = 1.004705065263158
$ws.Cells.Item(21, 8).Value2 = 1.015296312631579
$ws.Cells.Item(21, 9).Value2 = 0.9941186657894736
$ws.Cells.Item(21, 10).Value2 = 1.005055029473684
$ws.Cells.Item(21, 11).Value2 = 0.9912676147368424
$ws.Cells.Item(21, 12).Value2 = 0.984316445263158
$ws.Cells.Item(21, 13).Value2 = 0.984316445263158
$ws.Cells.Item(21, 14).Value2 = 0.984316445263158
$ws.Cells.Item(21, 15).Value2 = 0.9941186657894736
$ws.Cells.Item(21, 16).Value2 = 1.008822003947368
$ws.Cells.Item(21, 17).Value2 = 0.9995868476315789
$ws.Cells.Item(21, 18).Value2 = 1.000653484385965
$ws.Cells.Item(21, 19).Value2 = 1.007566345789473
$ws.Cells.Item(21, 20).Value2 = 1.000653484385965
$ws.Cells.Item(21, 21).Value2 = 1.001753870657895
$ws.Cells.Item(21, 22).Value2 = 0.9982663855789473
$ws.Cells.Item(21, 23).Value2 = 1.001550392631579
# Row 22
$ws.Cells.Item(22, 2).Value2 = 'Ring Perpendicular to TD'
$ws.Cells.Item(22, 3).Value2 = 1.003166876315789
$ws.Cells.Item(22, 4).Value2 = 1.001187588947368
$ws.Cells.Item(22, 5).Value2 = 0.9952496973684208
$ws.Cells.Item(22, 6).Value2 = 0.9952496973684208
$ws.Cells.Item(22, 7).Value2 = 0.9990499384210526
$ws.Cells.Item(22, 8).Value2 = 0.9969113226315791
$ws.Cells.Item(22, 9).Value2 = 1.001187588947368
$ws.Cells.Item(22, 10).Value2 = 0.9989792663157895
$ws.Cells.Item(22, 11).Value2 = 1.001763281052632
$ws.Cells.Item(22, 12).Value2 = 1.003166876315789
$ws.Cells.Item(22, 13).Value2 = 1.003166876315789
$ws.Cells.Item(22, 14).Value2 = 1.003166876315789
$ws.Cells.Item(22, 15).Value2 = 1.001187588947368
$ws.Cells.Item(22, 16).Value2 = 0.9982186431578944
$ws.Cells.Item(22, 17).Value2 = 1.000083427631579
$ws.Cells.Item(22, 18).Value2 = 0.999868054210526
$ws.Cells.Item(22, 19).Value2 = 0.9984721842105261
$ws.Cells.Item(22, 20).Value2 = 0.999868054210526
$ws.Cells.Item(22, 21).Value2 = 0.9996458572368419
$ws.Cells.Item(22, 22).Value2 = 1.000350061052631
$ws.Cells.Item(22, 23).Value2 = 0.9996869449999999
# Row 23
$ws.Cells.Item(23, 2).Value2 = 'OffsetFTD'
$ws.Cells.Item(23, 3).Value2 = 0.969532040033978
$ws.Cells.Item(23, 4).Value2 = 0.988574515860675
$ws.Cells.Item(23, 5).Value2 = 1.045701945667218
$ws.Cells.Item(23, 6).Value2 = 1.045701945667218
$ws.Cells.Item(23, 7).Value2 = 1.00914037468098
$ws.Cells.Item(23, 8).Value2 = 1.029715670191089
$ws.Cells.Item(23, 9).Value2 = 0.988574515860675
$ws.Cells.Item(23, 10).Value2 = 1.009820251921844
$ws.Cells.Item(23, 11).Value2 = 0.983035846683217
$ws.Cells.Item(23, 12).Value2 = 0.969532040033978
$ws.Cells.Item(23, 13).Value2 = 0.969532040033978
$ws.Cells.Item(23, 14).Value2 = 0.969532040033978
$ws.Cells.Item(23, 15).Value2 = 0.988574515860675
$ws.Cells.Item(23, 16).Value2 = 1.017138230763946
$ws.Cells.Item(23, 17).Value2 = 0.9991973838912598
$ws.Cells.Item(23, 18).Value2 = 1.001269500520624
$ws.Cells.Item(23, 19).Value2 = 1.014698904483246
$ws.Cells.Item(23, 20).Value2 = 1.001269500520624
$ws.Cells.Item(23, 21).Value2 = 1.003407188370929
$ws.Cells.Item(23, 22).Value2 = 0.9966321587035386
$ws.Cells.Item(23, 23).Value2 = 1.00301189511246
# Row 24
$ws.Cells.Item(24, 2).Value2 = 'OffsetATD'
$ws.Cells.Item(24, 3).Value2 = 0.9918791945001935
$ws.Cells.Item(24, 4).Value2 = 0.9969546944137678
$ws.Cells.Item(24, 5).Value2 = 1.012181208237609
$ws.Cells.Item(24, 6).Value2 = 1.012181208237609
$ws.Cells.Item(24, 7).Value2 = 1.002436247123406
$ws.Cells.Item(24, 8).Value2 = 1.007920291228564
$ws.Cells.Item(24, 9).Value2 = 0.9969546944137678
$ws.Cells.Item(24, 10).Value2 = 1.002617452885799
$ws.Cells.Item(24, 11).Value2 = 0.9954784415153435
$ws.Cells.Item(24, 12).Value2 = 0.9918791945001935
$ws.Cells.Item(24, 13).Value2 = 0.9918791945001935
$ws.Cells.Item(24, 14).Value2 = 0.9918791945001935
$ws.Cells.Item(24, 15).Value2 = 0.9969546944137678
$ws.Cells.Item(24, 16).Value2 = 1.004567951325689
$ws.Cells.Item(24, 17).Value2 = 0.9997860736497832
$ws.Cells.Item(24, 18).Value2 = 1.00033836571719
$ws.Cells.Item(24, 19).Value2 = 1.003917785179059
$ws.Cells.Item(24, 20).Value2 = 1.00033836571719
$ws.Cells.Item(24, 21).Value2 = 1.000908137509342
$ws.Cells.Item(24, 22).Value2 = 0.9991023489075127
$ws.Cells.Item(24, 23).Value2 = 1.000802778039807
# Row 25
$ws.Cells.Item(25, 2).Value2 = 'OffsetF45'
$ws.Cells.Item(25, 3).Value2 = 1.000320241021014
$ws.Cells.Item(25, 4).Value2 = 1.000120094216192
$ws.Cells.Item(25, 5).Value2 = 0.9995196357321918
$ws.Cells.Item(25, 6).Value2 = 0.9995196357321918
$ws.Cells.Item(25, 7).Value2 = 0.9999039259782335
$ws.Cells.Item(25, 8).Value2 = 0.9996876591446778
$ws.Cells.Item(25, 9).Value2 = 1.000120094216192
$ws.Cells.Item(25, 10).Value2 = 0.9998967808468452
$ws.Cells.Item(25, 11).Value2 = 1.0001783130487
$ws.Cells.Item(25, 12).Value2 = 1.000320241021014
$ws.Cells.Item(25, 13).Value2 = 1.000320241021014
$ws.Cells.Item(25, 14).Value2 = 1.000320241021014
$ws.Cells.Item(25, 15).Value2 = 1.000120094216192
$ws.Cells.Item(25, 16).Value2 = 0.9998198649741917
$ws.Cells.Item(25, 17).Value2 = 1.000008437531518
$ws.Cells.Item(25, 18).Value2 = 0.9999866569897993
$ws.Cells.Item(25, 19).Value2 = 0.9998455035984096
$ws.Cells.Item(25, 20).Value2 = 0.9999866569897993
$ws.Cells.Item(25, 21).Value2 = 0.9999641879540607
$ws.Cells.Item(25, 22).Value2 = 1.000035398567451
$ws.Cells.Item(25, 23).Value2 = 0.9999683430255056
# Row 26
$ws.Cells.Item(26, 2).Value2 = 'OffsetA45'
$ws.Cells.Item(26, 3).Value2 = 1.000086877963037
$ws.Cells.Item(26, 4).Value2 = 1.000032580109521
$ws.Cells.Item(26, 5).Value2 = 0.9998696851522243
$ws.Cells.Item(26, 6).Value2 = 0.9998696851522243
$ws.Cells.Item(26, 7).Value2 = 0.9999739352794422
$ws.Cells.Item(26, 8).Value2 = 0.9999152752583488
$ws.Cells.Item(26, 9).Value2 = 1.000032580109521
$ws.Cells.Item(26, 10).Value2 = 0.9999719971649959
$ws.Cells.Item(26, 11).Value2 = 1.000048373023327
$ws.Cells.Item(26, 12).Value2 = 1.000086877963037
$ws.Cells.Item(26, 13).Value2 = 1.000086877963037
$ws.Cells.Item(26, 14).Value2 = 1.000086877963037
$ws.Cells.Item(26, 15).Value2 = 1.000032580109521
$ws.Cells.Item(26, 16).Value2 = 0.9999511326308728
$ws.Cells.Item(26, 17).Value2 = 1.000002288637258
$ws.Cells.Item(26, 18).Value2 = 0.9999963810749274
$ws.Cells.Item(26, 19).Value2 = 0.9999580874755805
$ws.Cells.Item(26, 20).Value2 = 0.9999963810749272
$ws.Cells.Item(26, 21).Value2 = 0.9999902850974444
$ws.Cells.Item(26, 22).Value2 = 1.000009603670563
$ws.Cells.Item(26, 23).Value2 = 0.9999914130075521
# Row 27
$ws.Cells.Item(27, 2).Value2 = 'OffsetFRD'
$ws.Cells.Item(27, 3).Value2 = 1.023819690784401
$ws.Cells.Item(27, 4).Value2 = 1.008932394824936
$ws.Cells.Item(27, 5).Value2 = 0.9642704642877111
$ws.Cells.Item(27, 6).Value2 = 0.9642704642877111
$ws.Cells.Item(27, 7).Value2 = 0.9928540917298576
$ws.Cells.Item(27, 8).Value2 = 0.9767684495888942
$ws.Cells.Item(27, 9).Value2 = 1.008932394824936
$ws.Cells.Item(27, 10).Value2 = 0.9923225778369981
$ws.Cells.Item(27, 11).Value2 = 1.013262485136027
$ws.Cells.Item(27, 12).Value2 = 1.023819690784401
$ws.Cells.Item(27, 13).Value2 = 1.023819690784401
$ws.Cells.Item(27, 14).Value2 = 1.023819690784401
$ws.Cells.Item(27, 15).Value2 = 1.008932394824936
$ws.Cells.Item(27, 16).Value2 = 0.9866014295563237
$ws.Cells.Item(27, 17).Value2 = 1.000627486330967
$ws.Cells.Item(27, 18).Value2 = 0.9990075166323497
$ws.Cells.Item(27, 19).Value2 = 0.9885084789832151
$ws.Cells.Item(27, 20).Value2 = 0.9990075166323497
$ws.Cells.Item(27, 21).Value2 = 0.9973362819335118
$ws.Cells.Item(27, 22).Value2 = 1.00263296370369
$ws.Cells.Item(27, 23).Value2 = 0.9976453186267202
# Row 28
$ws.Cells.Item(28, 2).Value2 = 'OffsetARD'
$ws.Cells.Item(28, 3).Value2 = 1.006395918014905
$ws.Cells.Item(28, 4).Value2 = 1.002398478429737
$ws.Cells.Item(28, 5).Value2 = 0.9904061168911402
$ws.Cells.Item(28, 6).Value2 = 0.9904061168911402
$ws.Cells.Item(28, 7).Value2 = 0.9980812246031864
$ws.Cells.Item(28, 8).Value2 = 0.993761999521841
$ws.Cells.Item(28, 9).Value2 = 1.002398478429737
$ws.Cells.Item(28, 10).Value2 = 0.997938500438495
$ws.Cells.Item(28, 11).Value2 = 1.003561163670427
$ws.Cells.Item(28, 12).Value2 = 1.006395918014905
$ws.Cells.Item(28, 13).Value2 = 1.006395918014905
$ws.Cells.Item(28, 14).Value2 = 1.006395918014905
$ws.Cells.Item(28, 15).Value2 = 1.002398478429737
$ws.Cells.Item(28, 16).Value2 = 0.9964022976604388
$ws.Cells.Item(28, 17).Value2 = 1.000168489434116
$ws.Cells.Item(28, 18).Value2 = 0.9997335044452608
$ws.Cells.Item(28, 19).Value2 = 0.9969143652531242
$ws.Cells.Item(28, 20).Value2 = 0.9997335044452608
$ws.Cells.Item(28, 21).Value2 = 0.9992847534435694
$ws.Cells.Item(28, 22).Value2 = 1.000706986357837
$ws.Cells.Item(28, 23).Value2 = 0.9993677349999337
# Row 29
$ws.Cells.Item(29, 2).Value2 = 'Gaussian Quadrature'
$ws.Cells.Item(29, 3).Value2 = 1.000683081276603
$ws.Cells.Item(29, 4).Value2 = 1.000256170142117
$ws.Cells.Item(29, 5).Value2 = 0.9989753889113621
$ws.Cells.Item(29, 6).Value2 = 0.9989753889113621
$ws.Cells.Item(29, 7).Value2 = 0.9997950781118711
$ws.Cells.Item(29, 8).Value2 = 0.999333782218699
$ws.Cells.Item(29, 9).Value2 = 1.000256170142117
$ws.Cells.Item(29, 10).Value2 = 0.9997798326548842
$ws.Cells.Item(29, 11).Value2 = 1.00038033083073
$ws.Cells.Item(29, 12).Value2 = 1.000683081276603
$ws.Cells.Item(29, 13).Value2 = 1.000683081276603
$ws.Cells.Item(29, 14).Value2 = 1.000683081276603
$ws.Cells.Item(29, 15).Value2 = 1.000256170142117
$ws.Cells.Item(29, 16).Value2 = 0.9996157795267397
$ws.Cells.Item(29, 17).Value2 = 1.000018001398501
$ws.Cells.Item(29, 18).Value2 = 0.999971546776694
$ws.Cells.Item(29, 19).Value2 = 0.9996704639027879
$ws.Cells.Item(29, 20).Value2 = 0.999971546776694
$ws.Cells.Item(29, 21).Value2 = 0.9999236182462415
$ws.Cells.Item(29, 22).Value2 = 1.000075510852314
$ws.Cells.Item(29, 23).Value2 = 0.999932479286048
# Row 30
$ws.Cells.Item(30, 1).Value2 = 28
$ws.Cells.Item(30, 2).Value2 = 'Michael-CCHex'
$ws.Cells.Item(30, 3).Value2 = 0.9993108581496504
$ws.Cells.Item(30, 4).Value2 = 0.9997415824313458
$ws.Cells.Item(30, 5).Value2 = 1.001033708190077
$ws.Cells.Item(30, 6).Value2 = 1.001033708190077
$ws.Cells.Item(30, 7).Value2 = 1.000206742516375
$ws.Cells.Item(30, 8).Value2 = 1.000672125815247
$ws.Cells.Item(30, 9).Value2 = 0.9997415824313458
$ws.Cells.Item(30, 10).Value2 = 1.000222114366262
$ws.Cells.Item(30, 11).Value2 = 0.9996163043673572
$ws.Cells.Item(30, 12).Value2 = 0.9993108581496504
$ws.Cells.Item(30, 13).Value2 = 0.9993108581496504
$ws.Cells.Item(30, 14).Value2 = 0.9993108581496504
$ws.Cells.Item(30, 15).Value2 = 0.9997415824313458
$ws.Cells.Item(30, 16).Value2 = 1.000387645310711
$ws.Cells.Item(30, 17).Value2 = 0.9999818483988038
$ws.Cells.Item(30, 18).Value2 = 1.000028716257024
$ws.Cells.Item(30, 19).Value2 = 1.000332468329228
$ws.Cells.Item(30, 20).Value2 = 1.000028716257024
$ws.Cells.Item(30, 21).Value2 = 1.000077065784334
$ws.Cells.Item(30, 22).Value2 = 0.999923824257397
$ws.Cells.Item(30, 23).Value2 = 1.000068127283457
# Row 31
$ws.Cells.Item(31, 1).Value2 = 29
$ws.Cells.Item(31, 2).Value2 = 'Michael-SNHex'
$ws.Cells.Item(31, 3).Value2 = 0.9970030735112908
$ws.Cells.Item(31, 4).Value2 = 0.998876172065593
$ws.Cells.Item(31, 5).Value2 = 1.004495382678261
$ws.Cells.Item(31, 6).Value2 = 1.004495382678261
$ws.Cells.Item(31, 7).Value2 = 1.000899078994334
$ws.Cells.Item(31, 8).Value2 = 1.002922929995886
$ws.Cells.Item(31, 9).Value2 = 0.998876172065593
$ws.Cells.Item(31, 10).Value2 = 1.000965935724017
$ws.Cells.Item(31, 11).Value2 = 0.998331362652343
$ws.Cells.Item(31, 12).Value2 = 0.9970030735112908
$ws.Cells.Item(31, 13).Value2 = 0.9970030735112908
$ws.Cells.Item(31, 14).Value2 = 0.9970030735112908
$ws.Cells.Item(31, 15).Value2 = 0.998876172065593
$ws.Cells.Item(31, 16).Value2 = 1.001685777371927
$ws.Cells.Item(31, 17).Value2 = 0.999921053894805
$ws.Cells.Item(31, 18).Value2 = 1.000124876085048
$ws.Cells.Item(31, 19).Value2 = 1.001445830155957
$ws.Cells.Item(31, 20).Value2 = 1.000124876085048
$ws.Cells.Item(31, 21).Value2 = 1.00033514099479
$ws.Cells.Item(31, 22).Value2 = 0.9996687274980905
$ws.Cells.Item(31, 23).Value2 = 1.000296263460915
